$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 1.03
$ws.Range("O2").Value = 1.3
$ws.Range("AT2").Value = 2.62
$ws.Range("G3").Value = 2.35
$ws.Range("I3").Value = 3.4
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.03
$ws.Range("O3").Value = 1.3
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("X3").Value = 11
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 5.5
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 9.5
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 34
$ws.Range("AT3").Value = 2.62
$ws.Range("AV3").Value = 51
$ws.Range("AZ3").Value = 51
$ws.Range("BA3").Value = 81
$ws.Range("J4").Value = 2.37
$ws.Range("M4").Value = 1.03
$ws.Range("O4").Value = 1.38
$ws.Range("P4").Value = 2.65
$ws.Range("R5").Value = 1.5
$ws.Range("V5").Value = 1.73
$ws.Range("R6").Value = 1.62
$ws.Range("V6").Value = 1.73
$ws.Range("M7").Value = 1.04
$ws.Range("O7").Value = 1.22
$ws.Range("S7").Value = 1.33
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.29
$ws.Range("S8").Value = 1.37
$ws.Range("M9").Value = 1.05
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 1.93
$ws.Range("R9").Value = 1.93
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.75
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("AB9").Value = 41
$ws.Range("AG9").Value = 251
$ws.Range("AT9").Value = 2.75
$ws.Range("BB9").Value = 151
$ws.Range("P10").Value = 4.25
$ws.Range("S12").Value = 1.27
$ws.Range("S13").Value = 1.19
$ws.Range("N14").Value = 29
$ws.Range("S14").Value = 1.17
$ws.Range("S15").Value = 1.25
$ws.Range("Q16").Value = 1.6
$ws.Range("G17").Value = 1.5
$ws.Range("I18").Value = 2.87
$ws.Range("Q18").Value = 1.44
$ws.Range("U19").Value = 1.57
$ws.Range("G20").Value = 1.69
$ws.Range("U20").Value = 1.67
$ws.Range("S21").Value = 1.22
$ws.Range("U21").Value = 1.53
$ws.Range("V21").Value = 2.38
$ws.Range("J22").Value = 1.69
$ws.Range("K22").Value = 2.87
$ws.Range("S22").Value = 1.17
$ws.Range("U22").Value = 1.57
$ws.Range("S23").Value = 1.47
$ws.Range("U23").Value = 1.91
$ws.Range("V23").Value = 1.8
$ws.Range("G25").Value = 1.8
$ws.Range("Q25").Value = 1.95
$ws.Range("R25").Value = 1.9
$ws.Range("S25").Value = 1.37
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.95
$ws.Range("Q26").Value = 1.8
$ws.Range("U26").Value = 1.62
$ws.Range("R27").Value = 1.62
$ws.Range("V27").Value = 1.73
$ws.Range("U28").Value = 1.73
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("BD32").Value = 126
$ws.Range("R34").Value = 1.65
$ws.Range("J35").Value = 2.87
